# Generate Report for Handback
# Update the handoff/handback timestamp strings recorded on the
# "Overview", "zh-cn" and "de-de" sheets of the handback status report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-24 02:45:11"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-24 02:44:58"
$zhcn.Range("K3").Value = "2016-08-24 02:45:33"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-24 02:45:11"
$dede.Range("K3").Value = "2016-08-24 02:45:40"
